$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.055109983661837
$ws.Range("D2").Value = 1.056223900627459
$ws.Range("E2").Value = 1.061033423298675
$ws.Range("F2").Value = 1.069197292762999
$ws.Range("I2").Value = 1.039435015973638
$ws.Range("J2").Value = 1.060118561548566
$ws.Range("K2").Value = 1.058962199609864
$ws.Range("L2").Value = 1.063758583645706
$ws.Range("M2").Value = 1.071900443808422
$ws.Range("N2").Value = 1.061624051503264
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.05665208065883
$ws.Range("D3").Value = 1.057379343462143
$ws.Range("E3").Value = 1.062300387550669
$ws.Range("F3").Value = 1.070539574357999
$ws.Range("I3").Value = 1.039742735066519
$ws.Range("J3").Value = 1.061309342210756
$ws.Range("K3").Value = 1.059930368674904
$ws.Range("L3").Value = 1.064838949099772
$ws.Range("M3").Value = 1.073057544486944
$ws.Range("N3").Value = 1.062816523210578
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.057648547783036
$ws.Range("D4").Value = 1.058125668478533
$ws.Range("E4").Value = 1.063119854407849
$ws.Range("F4").Value = 1.071407415427298
$ws.Range("I4").Value = 1.039939911084884
$ws.Range("J4").Value = 1.062078092263129
$ws.Range("K4").Value = 1.060554937912298
$ws.Range("L4").Value = 1.065537116502639
$ws.Range("M4").Value = 1.073805017163598
$ws.Range("N4").Value = 1.063586364976203
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.058067141407053
$ws.Range("D5").Value = 1.058439111072918
$ws.Range("E5").Value = 1.063464280182639
$ws.Range("F5").Value = 1.071772091993664
$ws.Range("I5").Value = 1.04002234128499
$ws.Range("J5").Value = 1.062400857747331
$ws.Range("K5").Value = 1.060817055943632
$ws.Range("L5").Value = 1.065830414207714
$ws.Range("M5").Value = 1.074118960232562
$ws.Range("N5").Value = 1.063909588824414
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.058137406502522
$ws.Range("D6").Value = 1.058491721271269
$ws.Range("E6").Value = 1.063522106327987
$ws.Range("F6").Value = 1.071833313330127
$ws.Range("I6").Value = 1.040036154596367
$ws.Range("J6").Value = 1.062455027194153
$ws.Range("K6").Value = 1.060861040396437
$ws.Range("L6").Value = 1.065879647871627
$ws.Range("M6").Value = 1.074171655520027
$ws.Range("N6").Value = 1.063963835198064
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.057654142304195
$ws.Range("D7").Value = 1.058129857935997
$ws.Range("E7").Value = 1.063124456943564
$ws.Range("F7").Value = 1.07141228889103
$ws.Range("I7").Value = 1.039941014337033
$ws.Range("J7").Value = 1.062082406703309
$ws.Range("K7").Value = 1.060558442112725
$ws.Range("L7").Value = 1.065541036388622
$ws.Range("M7").Value = 1.073809213237627
$ws.Range("N7").Value = 1.063590685543382
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.055631431314546
$ws.Range("D8").Value = 1.05661466473425
$ws.Range("E8").Value = 1.061461673506043
$ws.Range("F8").Value = 1.069651072038054
$ws.Range("I8").Value = 1.0395394134663
$ws.Range("J8").Value = 1.060521360108059
$ws.Range("K8").Value = 1.059289793005454
$ws.Range("L8").Value = 1.064123887180407
$ws.Range("M8").Value = 1.07229175271401
$ws.Range("N8").Value = 1.062027422082913
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.052056306300109
$ws.Range("D9").Value = 1.053934361095671
$ws.Range("E9").Value = 1.058528815086144
$ws.Range("F9").Value = 1.066541962572081
$ws.Range("I9").Value = 1.038816822162628
$ws.Range("J9").Value = 1.057756837343677
$ws.Range("K9").Value = 1.057039523560742
$ws.Range("L9").Value = 1.061619610916271
$ws.Range("M9").Value = 1.069608024964468
$ws.Range("N9").Value = 1.05925897337911
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.049665122826519
$ws.Range("D10").Value = 1.052140256335079
$ws.Range("E10").Value = 1.056571412427316
$ws.Range("F10").Value = 1.064465131480744
$ws.Range("I10").Value = 1.038324964176235
$ws.Range("J10").Value = 1.055904227850795
$ws.Range("K10").Value = 1.055529161792107
$ws.Range("L10").Value = 1.059945086737904
$ws.Range("M10").Value = 1.06781201803012
$ws.Range("N10").Value = 1.057403732968251
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.048627766211048
$ws.Range("D11").Value = 1.051361609021401
$ws.Range("E11").Value = 1.055723262193132
$ws.Range("F11").Value = 1.063564799936879
$ws.Range("I11").Value = 1.038109558855991
$ws.Range("J11").Value = 1.05509967476707
$ws.Range("K11").Value = 1.054872683378277
$ws.Range("L11").Value = 1.059218755993161
$ws.Range("M11").Value = 1.067032636607379
$ws.Range("N11").Value = 1.056598037326863
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.048242142373272
$ws.Range("D12").Value = 1.051072110491333
$ws.Range("E12").Value = 1.055408128739113
$ws.Range("F12").Value = 1.063230212819784
$ws.Range("I12").Value = 1.0380291810057
$ws.Range("J12").Value = 1.054800466491975
$ws.Range("K12").Value = 1.054628459938331
$ws.Range("L12").Value = 1.058948771619306
$ws.Range("M12").Value = 1.066742878793895
$ws.Range("N12").Value = 1.056298404141698
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.048324873873143
$ws.Range("D13").Value = 1.051134221409556
$ws.Range("E13").Value = 1.055475730214191
$ws.Range("F13").Value = 1.063301990407963
$ws.Range("I13").Value = 1.038046438954883
$ws.Range("J13").Value = 1.054864664147079
$ws.Range("K13").Value = 1.05468086389127
$ws.Range("L13").Value = 1.059006692987269
$ws.Range("M13").Value = 1.066805044723909
$ws.Range("N13").Value = 1.056362692964836
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.048595896666957
$ws.Range("D14").Value = 1.051337684609768
$ws.Range("E14").Value = 1.055697215084998
$ws.Range("F14").Value = 1.063537146240532
$ws.Range("I14").Value = 1.038102922290542
$ws.Range("J14").Value = 1.055074949530041
$ws.Range("K14").Value = 1.054852503536834
$ws.Range("L14").Value = 1.059196442961975
$ws.Range("M14").Value = 1.067008690490311
$ws.Range("N14").Value = 1.056573276977161
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.04876284227569
$ws.Range("D15").Value = 1.051463008493509
$ws.Range("E15").Value = 1.055833666783045
$ws.Range("F15").Value = 1.06368201160778
$ws.Range("I15").Value = 1.038137674883016
$ws.Range("J15").Value = 1.055204465198605
$ws.Range("K15").Value = 1.054958206123523
$ws.Range("L15").Value = 1.05931332848979
$ws.Range("M15").Value = 1.067134128632753
$ws.Range("N15").Value = 1.056702976572831
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.049733926725833
$ws.Range("D16").Value = 1.05219189444204
$ws.Range("E16").Value = 1.056627688655714
$ws.Range("F16").Value = 1.06452486087472
$ws.Range("I16").Value = 1.038339208600797
$ws.Range("J16").Value = 1.055957573073383
$ws.Range("K16").Value = 1.055572677365259
$ws.Range("L16").Value = 1.059993264149626
$ws.Range("M16").Value = 1.067863706739294
$ws.Range("N16").Value = 1.057457153947174
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.050342531570598
$ws.Range("D17").Value = 1.052648622669146
$ws.Range("E17").Value = 1.057125598089295
$ws.Range("F17").Value = 1.065053272518164
$ws.Range("I17").Value = 1.038464974035965
$ws.Range("J17").Value = 1.056429340917086
$ws.Range("K17").Value = 1.055957450788487
$ws.Range("L17").Value = 1.060419431531435
$ws.Range("M17").Value = 1.06832089307266
$ws.Range("N17").Value = 1.057929591755327
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.050697332147125
$ws.Range("D18").Value = 1.052914852371527
$ws.Range("E18").Value = 1.057415964444862
$ws.Range("F18").Value = 1.065361385200783
$ws.Range("I18").Value = 1.038538096798516
$ws.Range("J18").Value = 1.056704287779364
$ws.Range("K18").Value = 1.056181643517034
$ws.Range("L18").Value = 1.060667887411073
$ws.Range("M18").Value = 1.068587398624079
$ws.Range("N18").Value = 1.058204929073686
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.050818278306702
$ws.Range("D19").Value = 1.053005600789765
$ws.Range("E19").Value = 1.057514962563039
$ws.Range("F19").Value = 1.065466426735597
$ws.Range("I19").Value = 1.038562990125869
$ws.Range("J19").Value = 1.056797999196029
$ws.Range("K19").Value = 1.056258047046763
$ws.Range("L19").Value = 1.060752584179921
$ws.Range("M19").Value = 1.068678242567503
$ws.Range("N19").Value = 1.058298773571313
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.050277253577185
$ws.Range("D20").Value = 1.052599637899284
$ws.Range("E20").Value = 1.057072182928384
$ws.Range("F20").Value = 1.06499658942971
$ws.Range("I20").Value = 1.038451504822412
$ws.Range("J20").Value = 1.056378748255815
$ws.Range("K20").Value = 1.055916193011007
$ws.Range("L20").Value = 1.060373720287816
$ws.Range("M20").Value = 1.06827185825871
$ws.Range("N20").Value = 1.057878927246674
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.04851609565866
$ws.Range("D21").Value = 1.051277777384332
$ws.Range("E21").Value = 1.055631995912138
$ws.Range("F21").Value = 1.06346790328284
$ws.Range("I21").Value = 1.038086299495582
$ws.Range("J21").Value = 1.05501303575773
$ws.Range("K21").Value = 1.054801970419877
$ws.Range("L21").Value = 1.05914057168536
$ws.Range("M21").Value = 1.066948729134289
$ws.Range("N21").Value = 1.056511275280192
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.047407025788812
$ws.Range("D22").Value = 1.050445081687378
$ws.Range("E22").Value = 1.054725954245154
$ws.Range("F22").Value = 1.062505806306683
$ws.Range("I22").Value = 1.037854557655327
$ws.Range("J22").Value = 1.054152263570355
$ws.Range("K22").Value = 1.054099224086835
$ws.Range("L22").Value = 1.058364124311884
$ws.Range("M22").Value = 1.066115313298913
$ws.Range("N22").Value = 1.05564928069758
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.047995134617439
$ws.Range("D23").Value = 1.050886661982674
$ws.Range("E23").Value = 1.055206316925836
$ws.Range("F23").Value = 1.063015924355589
$ws.Range("I23").Value = 1.037977610277258
$ws.Range("J23").Value = 1.054608775990912
$ws.Range("K23").Value = 1.054471972603898
$ws.Range("L23").Value = 1.058775841361026
$ws.Range("M23").Value = 1.06655726802654
$ws.Range("N23").Value = 1.056106441418136
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.050306750474734
$ws.Range("D24").Value = 1.052621772540473
$ws.Range("E24").Value = 1.057096319113111
$ws.Range("F24").Value = 1.065022202389562
$ws.Range("I24").Value = 1.038457591703667
$ws.Range("J24").Value = 1.056401609604421
$ws.Range("K24").Value = 1.055934836363771
$ws.Range("L24").Value = 1.060394375603333
$ws.Range("M24").Value = 1.068294015485954
$ws.Range("N24").Value = 1.057901821061018
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.052981895972063
$ws.Range("D25").Value = 1.054628538219917
$ws.Range("E25").Value = 1.059287390987335
$ws.Range("F25").Value = 1.067346439980018
$ws.Range("I25").Value = 1.039005407542354
$ws.Range("J25").Value = 1.058473198482194
$ws.Range("K25").Value = 1.057623047146733
$ws.Range("L25").Value = 1.062267890550889
$ws.Range("M25").Value = 1.070303019977433
$ws.Range("N25").Value = 1.059976351832612
